# Auto-generated Excel COM-interop edit script
# Applies per-cell numeric updates to match the target diff across 8 worksheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Source data has no formulas -
# every changed cell is a literal value, so we set .Value directly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 15626094
$ws.Range("I33").Value = 20833926
$ws.Range("K33").Value = 20833926
$ws.Range("M33").Value = -20833697

# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 3125
$ws.Range("I40").Value = 2250
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 2250
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -2075
$ws.Range("N40").Value = -4350

# Row 48: The Sting of Conscience / Sleeping Potion
$ws.Range("H48").Value = 2746
$ws.Range("J48").Value = 2746
$ws.Range("L48").Value = 8238
$ws.Range("N48").Value = -8822

# Row 56: Sleepless in Silvertear / Potent Sleeping Potion
$ws.Range("H56").Value = 2746
$ws.Range("J56").Value = 2746
$ws.Range("L56").Value = 8238
$ws.Range("N56").Value = -9306

# Row 69: Steeling the Knife, Steeling the Mind / Grade 1 Mind Dissolvent
$ws.Range("H69").Value = 8499.666999999999
$ws.Range("I69").Value = 7749.5
$ws.Range("K69").Value = 23248.5
$ws.Range("M69").Value = -22374.5

# Row 72: Surgical Substitution (L) / Grade 1 Mind Dissolvent
$ws.Range("H72").Value = 8499.666999999999
$ws.Range("I72").Value = 7749.5
$ws.Range("K72").Value = 69745.5
$ws.Range("M72").Value = -65377.5

# Row 80: Cleansing the Wicked Humours / Hallowed Water
$ws.Range("H80").Value = 1139.5555
$ws.Range("I80").Value = 642.7
$ws.Range("J80").Value = 1760.625
$ws.Range("K80").Value = 1928.1
$ws.Range("L80").Value = 5281.875
$ws.Range("M80").Value = -930.1000000000001
$ws.Range("N80").Value = -7277.875

# Row 83: Washing Away the Sins (L) / Hallowed Water
$ws.Range("H83").Value = 1139.5555
$ws.Range("I83").Value = 642.7
$ws.Range("J83").Value = 1760.625
$ws.Range("K83").Value = 5784.3
$ws.Range("L83").Value = 15845.625
$ws.Range("M83").Value = -792.3000000000002
$ws.Range("N83").Value = -25829.625

# Row 129: Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 1813.75
$ws.Range("I129").Value = 1038
$ws.Range("J129").Value = 2589.5
$ws.Range("K129").Value = 3114
$ws.Range("L129").Value = 7768.5
$ws.Range("M129").Value = 1886
$ws.Range("N129").Value = -17768.5

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 6387.875
$ws.Range("I132").Value = 5022.936
$ws.Range("K132").Value = 15068.808
$ws.Range("M132").Value = -12538.808

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 11176.473
$ws.Range("I137").Value = 2292.1
$ws.Range("K137").Value = 6876.299999999999
$ws.Range("M137").Value = -4326.299999999999


$ws = $wb.Worksheets.Item("ARM")

# Row 23: A Well-rounded Crew / Iron Hoplon
$ws.Range("H23").Value = 11999.667
$ws.Range("I23").Value = 11999.667
$ws.Range("K23").Value = 11999.667
$ws.Range("M23").Value = -11740.667


$ws = $wb.Worksheets.Item("BSM")

# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 14790.208
$ws.Range("I20").Value = 7367.607
$ws.Range("J20").Value = 25181.85
$ws.Range("K20").Value = 7367.607
$ws.Range("L20").Value = 25181.85
$ws.Range("M20").Value = -7120.607
$ws.Range("N20").Value = -25675.85

# Row 22: Riveting Run / Iron Rivets
$ws.Range("H22").Value = 1042.9
$ws.Range("I22").Value = 1174.875
$ws.Range("J22").Value = 515
$ws.Range("K22").Value = 1174.875
$ws.Range("L22").Value = 515
$ws.Range("M22").Value = -1001.875
$ws.Range("N22").Value = -861

# Row 25: Tools of the Trade / Iron Doming Hammer
$ws.Range("H25").Value = 1690.2727
$ws.Range("I25").Value = 510.33334
$ws.Range("K25").Value = 510.33334
$ws.Range("M25").Value = -275.33334

# Row 54: Get Me to the War on Time / Cobalt Doming Hammer
$ws.Range("H54").Value = 4294.3335
$ws.Range("I54").Value = 1441.5
$ws.Range("K54").Value = 1441.5
$ws.Range("M54").Value = -957.5

# Row 64: With Bearings Straight / Mythrite Nugget
$ws.Range("H64").Value = 2750
$ws.Range("I64").Value = 2000
$ws.Range("K64").Value = 2000
$ws.Range("M64").Value = -1775

# Row 67: Bearing the Brunt (L) / Mythrite Nugget
$ws.Range("H67").Value = 2750
$ws.Range("I67").Value = 2000
$ws.Range("K67").Value = 2000
$ws.Range("M67").Value = -1220

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 1143.1666
$ws.Range("I86").Value = 1164.55
$ws.Range("K86").Value = 1164.55
$ws.Range("M86").Value = -41.54999999999995

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 1143.1666
$ws.Range("I89").Value = 1164.55
$ws.Range("K89").Value = 5822.75
$ws.Range("M89").Value = -206.75


$ws = $wb.Worksheets.Item("CRP")

# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 3313
$ws.Range("I16").Value = 1004.26666
$ws.Range("J16").Value = 7160.8887
$ws.Range("K16").Value = 1004.26666
$ws.Range("L16").Value = 7160.8887
$ws.Range("M16").Value = -717.26666
$ws.Range("N16").Value = -7734.8887

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 6903.4517
$ws.Range("I99").Value = 2835.4
$ws.Range("J99").Value = 8840.619000000001
$ws.Range("K99").Value = 2835.4
$ws.Range("L99").Value = 8840.619000000001
$ws.Range("M99").Value = -1337.4
$ws.Range("N99").Value = -11836.619

# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 3313
$ws.Range("I113").Value = 1004.26666
$ws.Range("J113").Value = 7160.8887
$ws.Range("K113").Value = 1004.26666
$ws.Range("L113").Value = 7160.8887
$ws.Range("M113").Value = 1165.73334
$ws.Range("N113").Value = -11500.8887

# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 4776.3076
$ws.Range("I122").Value = 2337.625
$ws.Range("J122").Value = 8678.200000000001
$ws.Range("K122").Value = 7012.875
$ws.Range("L122").Value = 26034.6
$ws.Range("M122").Value = -4562.875
$ws.Range("N122").Value = -30934.6

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 6903.4517
$ws.Range("I126").Value = 2835.4
$ws.Range("J126").Value = 8840.619000000001
$ws.Range("K126").Value = 8506.200000000001
$ws.Range("L126").Value = 26521.857
$ws.Range("M126").Value = -6036.200000000001
$ws.Range("N126").Value = -31461.857


$ws = $wb.Worksheets.Item("CUL")

# Row 2: Pork Is a Salty Food / Table Salt
$ws.Range("H2").Value = 253.69565
$ws.Range("I2").Value = 363.6
$ws.Range("J2").Value = 47.625
$ws.Range("K2").Value = 2181.6
$ws.Range("L2").Value = 285.75
$ws.Range("M2").Value = -2068.6
$ws.Range("N2").Value = -511.75

# Row 38: Pretty as a Picture / Dark Vinegar
$ws.Range("H38").Value = 112.36364
$ws.Range("I38").Value = 20.4
$ws.Range("J38").Value = 189
$ws.Range("K38").Value = 61.2
$ws.Range("L38").Value = 567
$ws.Range("M38").Value = 285.8
$ws.Range("N38").Value = -1261

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 1409.2245
$ws.Range("J131").Value = 1449.978
$ws.Range("L131").Value = 4349.934
$ws.Range("N131").Value = -14429.934


$ws = $wb.Worksheets.Item("GSM")

# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 13719.429
$ws.Range("I70").Value = 17267.625
$ws.Range("J70").Value = 11535.923
$ws.Range("K70").Value = 17267.625
$ws.Range("L70").Value = 11535.923
$ws.Range("M70").Value = -16997.625
$ws.Range("N70").Value = -12075.923

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 13719.429
$ws.Range("I73").Value = 17267.625
$ws.Range("J73").Value = 11535.923
$ws.Range("K73").Value = 17267.625
$ws.Range("L73").Value = 11535.923
$ws.Range("M73").Value = -16331.625
$ws.Range("N73").Value = -13407.923

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 3522.2222
$ws.Range("I122").Value = 3449.75
$ws.Range("K122").Value = 10349.25
$ws.Range("M122").Value = -7899.25


$ws = $wb.Worksheets.Item("LTW")

# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value = 2691.9333
$ws.Range("I55").Value = 1452.6154
$ws.Range("K55").Value = 1452.6154
$ws.Range("M55").Value = -1279.6154

# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 15976.571
$ws.Range("J68").Value = 33986
$ws.Range("L68").Value = 33986
$ws.Range("N68").Value = -35484

# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 15976.571
$ws.Range("J71").Value = 33986
$ws.Range("L71").Value = 169930
$ws.Range("N71").Value = -177418

# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 5976.278
$ws.Range("I82").Value = 5225.1816
$ws.Range("J82").Value = 7156.5713
$ws.Range("K82").Value = 5225.1816
$ws.Range("L82").Value = 7156.5713
$ws.Range("M82").Value = -4864.1816
$ws.Range("N82").Value = -7878.5713

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 5976.278
$ws.Range("I85").Value = 5225.1816
$ws.Range("J85").Value = 7156.5713
$ws.Range("K85").Value = 5225.1816
$ws.Range("L85").Value = 7156.5713
$ws.Range("M85").Value = -3977.1816
$ws.Range("N85").Value = -9652.5713

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 7860.846
$ws.Range("I122").Value = 5673
$ws.Range("K122").Value = 17019
$ws.Range("M122").Value = -14569

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 2510265.8
$ws.Range("I132").Value = 4134.7
$ws.Range("J132").Value = 6687151
$ws.Range("K132").Value = 12404.1
$ws.Range("L132").Value = 20061453
$ws.Range("M132").Value = -9874.099999999999
$ws.Range("N132").Value = -20066513


$ws = $wb.Worksheets.Item("WVR")

# Row 2: The Unmentionables / Hempen Underpants
$ws.Range("H2").Value = 2092.5334
$ws.Range("J2").Value = 33
$ws.Range("L2").Value = 33
$ws.Range("N2").Value = -257

# Row 17: Making Gloves Out of Nothing at All / Hempen Bracers
$ws.Range("H17").Value = 4000
$ws.Range("I17").Value = 4000
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 4000
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -3828
$ws.Range("N17").ClearContents()

# Row 70: An Account of My Boots / Holy Rainbow Shoes
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73: Soot in My Hair and Scars on My Feet (L) / Holy Rainbow Shoes
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

